$d = $word.ActiveDocument

# --------------------------------------------------------------------
# 1. "Tuesday, March 10" -> "Tuesday, March 11"
#    (the replace merges the following superscript "th"/" " runs as a
#    side effect of this engine's run re-serialization; restore the
#    original 2-run split with a harmless Bold-toggle on just "th",
#    which forces Word to re-split the run without adding a spurious
#    xml:space="preserve" to a value that doesn't need it)
# --------------------------------------------------------------------
$idx = $d.Content.Text.IndexOf("Tuesday, March 10")
$len = "Tuesday, March 10".Length
$d.Range($idx, $idx + $len).Text = "Tuesday, March 11"

$txt = $d.Content.Text
$thIdx = $txt.IndexOf("th ", $idx)
$rTh = $d.Range($thIdx, $thIdx + 2)
$rTh.Bold = 1
$rTh.Bold = 0

# --------------------------------------------------------------------
# 2. Split "Begin data flow model for all components (UI and backend)"
#    into two runs, with the _GoBack bookmark moved in between them
#    (Word re-homes the hidden _GoBack bookmark at the most recent
#    edit location -- adding a bookmark of the same name relocates it)
# --------------------------------------------------------------------
$txt = $d.Content.Text
$full = "Begin data flow model for all components (UI and backend)"
$idx = $txt.IndexOf($full)
$splitPoint = $idx + "Begin data flow model for all components ".Length
$d.Bookmarks.Add("_GoBack", $d.Range($splitPoint, $splitPoint))

# --------------------------------------------------------------------
# 3. Merge "Continue working on pseu" + "docode for Library and
#    Definitions" (previously split by the old _GoBack bookmark
#    location) back into a single run now that the bookmark moved.
# --------------------------------------------------------------------
$d.Content.Find.Execute("Continue working on pseudocode for Library and Definitions",
                         $true, $false, $false, $false, $false, $true, 1, $false,
                         "Continue working on pseudocode for Library and Definitions", 2)
